# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Numeric-looking text values (e.g. "317.04") are prefixed with a leading
# apostrophe so Excel keeps storing them as text (quote-prefixed), matching
# the original inline-string cell content instead of converting to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.210.12'
$ws.Range('E2').Value = '  -1.02%  '

$ws.Range('D3').Value = '2.302.19'
$ws.Range('E3').Value = '  -2.23%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = '''317.04'
$ws.Range('E5').Value = '  -0.23%  '

$ws.Range('D6').Value = '''104.75'
$ws.Range('E6').Value = '  -2.40%  '

$ws.Range('D7').Value = '''0.628'
$ws.Range('E7').Value = '  -1.34%  '

$ws.Range('E8').Value = '  +0.12%  '

$ws.Range('D9').Value = '''0.610'
$ws.Range('E9').Value = '  -1.87%  '

$ws.Range('D10').Value = '''39.71'
$ws.Range('E10').Value = '  -4.35%  '

$ws.Range('D11').Value = '''0.0909'
$ws.Range('E11').Value = '  -1.88%  '

$ws.Range('D12').Value = '''8.44'
$ws.Range('E12').Value = '  -0.07%  '

$ws.Range('D13').Value = '''0.106'
$ws.Range('E13').Value = '  +0.39%  '

$ws.Range('D14').Value = '''0.977'
$ws.Range('E14').Value = '  -1.98%  '

$ws.Range('D15').Value = '''15.45'

$ws.Range('D16').Value = '2.651.66'

$ws.Range('D17').Value = '2.300.16'
$ws.Range('E17').Value = '  -1.66%  '

$ws.Range('D18').Value = '42.133.80'
$ws.Range('E18').Value = '  -1.24%  '

$ws.Range('D19').Value = '''7.76'
$ws.Range('E19').Value = '  -0.79%  '

$ws.Range('E20').Value = '  -0.43%  '

$ws.Range('D21').Value = '''286.23'
$ws.Range('E21').Value = '  +11.40%  '

$ws.Range('D22').Value = '''73.67'
$ws.Range('E22').Value = '  -3.91%  '

$ws.Range('E23').Value = '  -0.92%  '

$ws.Range('E24').Value = '  -0.82%  '

$ws.Range('D25').Value = '''9.96'
$ws.Range('E25').Value = '  +5.68%  '

$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.42%  '

$ws.Range('D27').Value = '''10.97'
$ws.Range('E27').Value = '  -3.68%  '

$ws.Range('D28').Value = '''23.42'
$ws.Range('E28').Value = '  +2.16%  '

$ws.Range('D29').Value = '''2.27'
$ws.Range('E29').Value = '  +1.79%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''164.99'
$ws.Range('E30').Value = '  -5.69%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '''35.50'
$ws.Range('E31').Value = '  -2.56%  '

$ws.Range('D32').Value = '''0.0883'
$ws.Range('E32').Value = '  -1.04%  '

$ws.Range('E33').Value = '  -0.87%  '

$ws.Range('D34').Value = '''5.92'
$ws.Range('E34').Value = '  -3.22%  '

$ws.Range('E35').Value = '  +1.52%  '

$ws.Range('D36').Value = '''0.117'
$ws.Range('E36').Value = '  -7.44%  '

$ws.Range('D37').Value = '''4.64'
$ws.Range('E37').Value = '  +0.80%  '

$ws.Range('D38').Value = '''2.94'
$ws.Range('E38').Value = '  +10.15%  '

$ws.Range('D39').Value = '''0.0353'
$ws.Range('E39').Value = '  -2.33%  '

$ws.Range('D40').Value = '''3.63'
$ws.Range('E40').Value = '  -4.15%  '

$ws.Range('D41').Value = '''103.09'
$ws.Range('E41').Value = '  +21.01%  '

$ws.Range('E42').Value = '  +1.88%  '

$ws.Range('D43').Value = '''71.18'
$ws.Range('E43').Value = '  -0.59%  '

$ws.Range('D44').Value = '''0.227'
$ws.Range('E44').Value = '  -5.28%  '

$ws.Range('E45').Value = '  +0.36%  '

$ws.Range('D46').Value = '''116.13'
$ws.Range('E46').Value = '  +1.66%  '

$ws.Range('D47').Value = '''12.12'
$ws.Range('E47').Value = '  +0.87%  '

$ws.Range('D48').Value = '''78.57'
$ws.Range('E48').Value = '  +4.23%  '

$ws.Range('E49').Value = '  +0.80%  '

$ws.Range('D50').Value = '''5.35'
$ws.Range('E50').Value = '  -2.54%  '

$ws.Range('D51').Value = '''1.28'
$ws.Range('E51').Value = '  +1.98%  '
